# Re-create the "open in a different Excel build, rename sheet, move the
# selection, nudge column A's width" edit described by the commit
# ("unify the conception of DataNode, DataTable, Entity.").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# The user's cursor ends up on C24 when the file was saved.
$ws.Range("C24").Select() | Out-Null

# Column A's width was nudged slightly (character-width units).
$ws.Columns(1).ColumnWidth = 23.4
